$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("postsArabic")

$ws.Cells.Item(151, 4).Value = "حسابات"
$ws.Cells.Item(152, 4).Value = "قروض"
$ws.Cells.Item(153, 4).Value = "معلومة"
$ws.Cells.Item(154, 4).Value = "قروض"
$ws.Cells.Item(155, 4).Value = "معلومة"
$ws.Cells.Item(156, 4).Value = "بطاقات"
$ws.Cells.Item(157, 4).Value = "قروض"
$ws.Cells.Item(158, 4).Value = "مساعدة"
$ws.Cells.Item(159, 4).Value = "معلومة"
$ws.Cells.Item(160, 4).Value = "معلومة"
$ws.Cells.Item(161, 4).Value = "تحويلات"
$ws.Cells.Item(162, 4).Value = "قروض"
$ws.Cells.Item(163, 4).Value = "مساعدة"
$ws.Cells.Item(164, 4).Value = "معلومة"
$ws.Cells.Item(165, 4).Value = "مساعدة"
$ws.Cells.Item(166, 4).Value = "مساعدة"
$ws.Cells.Item(167, 4).Value = "معلومة"
$ws.Cells.Item(168, 4).Value = "قروض"
$ws.Cells.Item(169, 4).Value = "قروض"
$ws.Cells.Item(170, 4).Value = "قروض"
$ws.Cells.Item(171, 4).Value = "مساعدة"
$ws.Cells.Item(172, 4).Value = "معلومة"
$ws.Cells.Item(173, 4).Value = "قروض"
$ws.Cells.Item(174, 4).Value = "بطاقات"
$ws.Cells.Item(175, 4).Value = "قروض"
$ws.Cells.Item(176, 4).Value = "قروض"
$ws.Cells.Item(177, 4).Value = "قروض"
$ws.Cells.Item(178, 4).Value = "معلومة"
$ws.Cells.Item(179, 4).Value = "حسابات"
$ws.Cells.Item(180, 4).Value = "معلومة"
$ws.Cells.Item(181, 4).Value = "قروض"
$ws.Cells.Item(182, 4).Value = "قروض"
$ws.Cells.Item(183, 4).Value = "قروض"
$ws.Cells.Item(184, 4).Value = "عملات"
$ws.Cells.Item(185, 4).Value = "قروض"
$ws.Cells.Item(186, 4).Value = "قروض"
$ws.Cells.Item(187, 4).Value = "عملات"
$ws.Cells.Item(188, 4).Value = "عملات"
$ws.Cells.Item(189, 4).Value = "حسابات"
$ws.Cells.Item(190, 4).Value = "قروض"
$ws.Cells.Item(191, 4).Value = "قروض"
$ws.Cells.Item(192, 4).Value = "ودائع"
$ws.Cells.Item(193, 4).Value = "تحويلات"
$ws.Cells.Item(194, 4).Value = "مساعدة"
$ws.Cells.Item(195, 4).Value = "حسابات"
$ws.Cells.Item(196, 4).Value = "قروض"
$ws.Cells.Item(197, 4).Value = "بطاقات"
$ws.Cells.Item(198, 4).Value = "مساعدة"
$ws.Cells.Item(199, 4).Value = "بطاقات"
$ws.Cells.Item(200, 4).Value = "مساعدة"

$ws.Range("E195").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 190
